$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 14 data rows (rows 2-15), shifting rows 16-22 up to become rows 2-8
$ws.Range("A2:A15").EntireRow.Delete() | Out-Null

# Append 13 new rows of data (rows 9-21) for the struggle class
$timestamps = @(700,800,900,1000,1100,1200,1300,1400,1500,1600,1700,1800,1900)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $r = 9 + $i
    $ws.Cells.Item($r, 1).Value2 = $timestamps[$i]
    $ws.Cells.Item($r, 2).Value2 = "struggle"
}

# Set the C:H numeric values for all data rows (2-21) to match the new dataset
$values = @{
    "C2" = -0.2141320705413827
    "D2" = -2.326400130987168
    "E2" = -2.468529105186464
    "F2" = 2.08137059211731
    "G2" = -2.974608421325684
    "H2" = 1.146135926246643
    "C3" = 2.301481246948242
    "D3" = -1.784020185470581
    "E3" = 0.5732603073120117
    "F3" = -0.7811439037322998
    "G3" = 1.68369734287262
    "H3" = -1.490358471870422
    "C4" = 0.7297354936599727
    "D4" = -2.234472751617432
    "E4" = -4.658630132675173
    "F4" = -0.4285219609737396
    "G4" = 1.809230089187622
    "H4" = -1.255480766296387
    "C5" = 1.63190019130707
    "D5" = -2.561783850193024
    "E5" = -1.591027021408075
    "F5" = -1.154382586479187
    "G5" = 3.296534299850464
    "H5" = -1.97553825378418
    "C6" = 0.6863539814949051
    "D6" = -1.503557689487935
    "E6" = 1.746999144554136
    "F6" = -1.183093309402466
    "G6" = 0.8458956480026245
    "H6" = -1.687667965888977
    "C7" = 0.2367095947265618
    "D7" = -1.311740666627883
    "E7" = -0.07715380191802312
    "F7" = -0.4327980279922485
    "G7" = 2.791501522064209
    "H7" = -2.889545440673828
    "C8" = 1.296695142984393
    "D8" = -3.332512527704244
    "E8" = -2.418631196022031
    "F8" = -2.166738986968994
    "G8" = 1.760666370391846
    "H8" = -2.120465993881226
    "C9" = 1.87941366434097
    "D9" = -4.668229699134821
    "E9" = -1.045372545719146
    "F9" = 1.305724501609802
    "G9" = 0.0389426611363887
    "H9" = -1.967597007751465
    "C10" = 0.3820920586585984
    "D10" = -1.311929136514663
    "E10" = 0.04638075828551969
    "F10" = 2.605340242385864
    "G10" = -0.3645338416099548
    "H10" = -0.5674937963485718
    "C11" = 2.75743055343628
    "D11" = -3.144901037216187
    "E11" = 4.137303829193115
    "F11" = -0.7612907886505127
    "G11" = 2.151620149612427
    "H11" = 0.1557706445455551
    "C12" = 4.460695505142212
    "D12" = 1.830066174268722
    "E12" = -0.1949661374092102
    "F12" = 0.5158756971359253
    "G12" = -2.014939069747925
    "H12" = 1.058934926986694
    "C13" = -3.40113162994386
    "D13" = 3.317261695861818
    "E13" = 1.794482350349431
    "F13" = -0.7177666425704956
    "G13" = 0.8080220222473145
    "H13" = 1.101084589958191
    "C14" = 0.008035421371476836
    "D14" = -0.4300747811794388
    "E14" = 2.762799173593522
    "F14" = 0.3738495409488678
    "G14" = -2.970790386199951
    "H14" = 2.029447078704834
    "C15" = 1.565377473831186
    "D15" = -1.947239398956303
    "E15" = 1.930309116840367
    "F15" = 1.252579212188721
    "G15" = -2.954449653625488
    "H15" = 2.593122959136963
    "C16" = -0.791193664073943
    "D16" = -1.430967807769773
    "E16" = 0.315328881144527
    "F16" = 1.917964100837708
    "G16" = -2.736523628234864
    "H16" = 1.525177836418152
    "C17" = 0.01700598001480103
    "D17" = -0.8887928128242493
    "E17" = 0.1178494691848755
    "F17" = -0.9870055317878724
    "G17" = -1.847256541252136
    "H17" = 0.9668469429016112
    "C18" = 3.75538071990013
    "D18" = 1.982138156890872
    "E18" = 1.528477013111122
    "F18" = 1.084591269493103
    "G18" = -5.03795862197876
    "H18" = 1.1534663438797
    "C19" = 4.293412685394287
    "D19" = -2.266220092773437
    "E19" = -9.760974884033203
    "F19" = -1.704924941062927
    "G19" = -0.2585487067699432
    "H19" = -1.87276017665863
    "C20" = 2.267539381980895
    "D20" = -2.359471559524536
    "E20" = -1.19436234235763
    "F20" = 0.7817547917366028
    "G20" = 1.986839175224304
    "H20" = -2.802802562713623
    "C21" = -1.645043730735778
    "D21" = -1.020936071872711
    "E21" = 0.8753915429115292
    "F21" = 0.7434229850769043
    "G21" = 0.6293439269065857
    "H21" = -2.307848930358887
}
foreach ($key in $values.Keys) {
    $ws.Range($key).Value2 = $values[$key]
}